$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the data rows for years 2000-2009 (original rows 2-11).
# This shifts the 2010-2018 rows (originally 12-20) up to become rows 2-10.
$ws.Rows("2:11").Delete()
